$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text looks like a plain number need to be forced to
# text format first, otherwise Excel auto-converts "1.00" -> 1, "0.0557" -> 5.57E-02, etc.
# We set NumberFormat = "@" (Text) before assigning, then restore the cell
# style to "Normal" afterwards so no visible formatting change remains.
$textCells = @("D5", "D6", "D8", "D11", "D20", "D21", "D22", "D24", "D25", "D26", "D31", "D32", "D33", "D38", "D40", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = "491.04"
$ws.Range("D6").Value = "151.56"
$ws.Range("D8").Value = "0.516"
$ws.Range("D11").Value = "5.72"
$ws.Range("D20").Value = "10.20"
$ws.Range("D21").Value = "321.06"
$ws.Range("D22").Value = "0.998"
$ws.Range("D24").Value = "58.16"
$ws.Range("D25").Value = "0.406"
$ws.Range("D26").Value = "1.00"
$ws.Range("D31").Value = "0.999"
$ws.Range("D32").Value = "151.15"
$ws.Range("D33").Value = "18.32"
$ws.Range("D38").Value = "3.77"
$ws.Range("D40").Value = "34.20"
$ws.Range("D42").Value = "0.0557"
$ws.Range("D43").Value = "0.996"
$ws.Range("D44").Value = "0.611"
$ws.Range("D45").Value = "0.0945"
$ws.Range("D46").Value = "4.81"
$ws.Range("D47").Value = "262.64"
$ws.Range("D48").Value = "10.20"
$ws.Range("D49").Value = "0.0229"

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = "57.006.99"
$ws.Range("E2").Value = "  +4.19%  "
$ws.Range("D3").Value = "2.479.08"
$ws.Range("E3").Value = "  +2.25%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +2.89%  "
$ws.Range("E6").Value = "  +9.24%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  +3.57%  "
$ws.Range("D9").Value = "2.492.10"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("E10").Value = "  +4.55%  "
$ws.Range("E11").Value = "  +4.71%  "
$ws.Range("E12").Value = "  +4.62%  "
$ws.Range("E13").Value = "  +1.60%  "
$ws.Range("D14").Value = "2.915.70"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").Value = "57.226.42"
$ws.Range("E15").Value = "  +4.26%  "
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("E17").Value = "  +2.57%  "
$ws.Range("D18").Value = "2.497.89"
$ws.Range("E18").Value = "  +1.94%  "
$ws.Range("E19").Value = "  +5.35%  "
$ws.Range("E20").Value = "  +4.06%  "
$ws.Range("E21").Value = "  +3.09%  "
$ws.Range("E22").Value = "  +0.37%  "
$ws.Range("E23").Value = "  +4.67%  "
$ws.Range("E24").Value = "  +2.03%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("E25").Value = "  +1.65%  "
$ws.Range("B26").Value = "Binance-PegBSC-USD"
$ws.Range("C26").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +2.30%  "
$ws.Range("D28").Value = "2.607.91"
$ws.Range("E29").Value = "  +2.95%  "
$ws.Range("D30").Value = "0.0₃0810"
$ws.Range("E30").Value = "  +6.06%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  +1.93%  "
$ws.Range("E33").Value = "  +2.71%  "
$ws.Range("E34").Value = "  +3.47%  "
$ws.Range("E35").Value = "  +1.74%  "
$ws.Range("E36").Value = "  +6.29%  "
$ws.Range("E37").Value = "  +3.60%  "
$ws.Range("E38").Value = "  +5.40%  "
$ws.Range("E39").Value = "  +9.06%  "
$ws.Range("E40").Value = "  +2.50%  "
$ws.Range("E41").Value = "  +3.42%  "
$ws.Range("E42").Value = "  +3.23%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("E43").Value = "  +0.27%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("E45").Value = "  +6.62%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E46").Value = "  +3.95%  "
$ws.Range("E47").Value = "  +3.30%  "
$ws.Range("E48").Value = "  +0.92%  "
$ws.Range("E49").Value = "  +3.77%  "
$ws.Range("E50").Value = "  +4.74%  "
$ws.Range("D51").Value = "1.870.82"
$ws.Range("E51").Value = "  -2.32%  "

